$wb = $excel.ActiveWorkbook

# ALC row 32
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1000
$ws.Range("J32").Value = 1000
$ws.Range("L32").Value = 1000
$ws.Range("N32").Value = -1652

# ALC row 98
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 0
$ws.Range("I98").Value = 0
$ws.Range("K98").Value = 0
$ws.Range("M98").ClearContents()

# ALC row 112
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 2562.2083
$ws.Range("J112").Value = 2567.5652
$ws.Range("L112").Value = 7702.6956
$ws.Range("N112").Value = -9918.695599999999

# ALC row 122
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()

# ARM row 97
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 2366.1667
$ws.Range("I97").Value = 2039.6
$ws.Range("J97").Value = 3999
$ws.Range("K97").Value = 2039.6
$ws.Range("L97").Value = 3999
$ws.Range("M97").Value = -1543.6
$ws.Range("N97").Value = -4991

# BSM row 20
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 6523.375
$ws.Range("I20").Value = 1872
$ws.Range("J20").Value = 11174.75
$ws.Range("K20").Value = 1872
$ws.Range("L20").Value = 11174.75
$ws.Range("M20").Value = -1625
$ws.Range("N20").Value = -11668.75

# BSM row 107
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 638.8
$ws.Range("I107").Value = 638.8
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 638.8
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 1281.2
$ws.Range("N107").ClearContents()

# CRP row 22
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 397.5
$ws.Range("I22").Value = 397.5
$ws.Range("K22").Value = 397.5
$ws.Range("M22").Value = -47.5

# CRP row 23
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H23").Value = 50000
$ws.Range("I23").Value = 50000
$ws.Range("K23").Value = 50000
$ws.Range("M23").Value = -49760

# CRP row 27
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H27").Value = 50000
$ws.Range("I27").Value = 50000
$ws.Range("K27").Value = 50000
$ws.Range("M27").Value = -49808

# CRP row 86
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 19503.143
$ws.Range("I86").Value = 34840.668
$ws.Range("K86").Value = 34840.668
$ws.Range("M86").Value = -33717.668

# CRP row 89
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H89").Value = 19503.143
$ws.Range("I89").Value = 34840.668
$ws.Range("K89").Value = 174203.34
$ws.Range("M89").Value = -168587.34

# CUL row 22
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 1499.8334
$ws.Range("I22").Value = 1499
$ws.Range("K22").Value = 4497
$ws.Range("M22").Value = -4328

# CUL row 27
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H27").Value = 1499.8334
$ws.Range("I27").Value = 1499
$ws.Range("K27").Value = 4497
$ws.Range("M27").Value = -4395

# CUL row 132
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 2499.6667
$ws.Range("I132").Value = 2499.6667
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 22497.0003
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -19967.0003
$ws.Range("N132").ClearContents()

# GSM row 13
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H13").Value = 6499.5
$ws.Range("J13").Value = 6499.5
$ws.Range("L13").Value = 6499.5
$ws.Range("N13").Value = -6777.5

# GSM row 28
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H28").Value = 9000
$ws.Range("J28").Value = 9000
$ws.Range("L28").Value = 9000
$ws.Range("N28").Value = -9384

# GSM row 107
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 2900
$ws.Range("I107").Value = 3000
$ws.Range("J107").Value = 2800
$ws.Range("K107").Value = 3000
$ws.Range("L107").Value = 2800
$ws.Range("M107").Value = -1080
$ws.Range("N107").Value = -6640

# GSM row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2933.111
$ws.Range("I122").Value = 2599.7144
$ws.Range("J122").Value = 4100
$ws.Range("K122").Value = 7799.1432
$ws.Range("L122").Value = 12300
$ws.Range("M122").Value = -5349.1432
$ws.Range("N122").Value = -17200

# GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4381.5
$ws.Range("I132").Value = 3245.5
$ws.Range("K132").Value = 9736.5
$ws.Range("M132").Value = -7206.5

# LTW row 16
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1174.4
$ws.Range("I16").Value = 1219.5
$ws.Range("K16").Value = 1219.5
$ws.Range("M16").Value = -1049.5

# LTW row 22
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 6817.385
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()

# LTW row 24
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H24").Value = 23737.5
$ws.Range("J24").Value = 50000
$ws.Range("L24").Value = 50000
$ws.Range("N24").Value = -50686

# LTW row 27
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 6817.385
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").ClearContents()

# LTW row 43
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()

# LTW row 55
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 2597.8
$ws.Range("I55").Value = 2997.5
$ws.Range("J55").Value = 999
$ws.Range("K55").Value = 2997.5
$ws.Range("L55").Value = 999
$ws.Range("M55").Value = -2824.5
$ws.Range("N55").Value = -1345

# LTW row 76
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H76").Value = 52816.43
$ws.Range("J76").Value = 52816.43
$ws.Range("L76").Value = 52816.43
$ws.Range("N76").Value = -53492.43

# LTW row 79
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H79").Value = 52816.43
$ws.Range("J79").Value = 52816.43
$ws.Range("L79").Value = 52816.43
$ws.Range("N79").Value = -55156.43

# LTW row 106
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H106").Value = 52132.332
$ws.Range("J106").Value = 52132.332
$ws.Range("L106").Value = 52132.332
$ws.Range("N106").Value = -54656.332

# WVR row 34
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H34").Value = 20000
$ws.Range("I34").Value = 20000
$ws.Range("K34").Value = 20000
$ws.Range("M34").Value = -19797

# WVR row 62
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 7958.1665
$ws.Range("I62").Value = 2449.5
$ws.Range("J62").Value = 10712.5
$ws.Range("K62").Value = 2449.5
$ws.Range("L62").Value = 10712.5
$ws.Range("M62").Value = -1825.5
$ws.Range("N62").Value = -11960.5

# WVR row 63
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 32830.8
$ws.Range("J63").Value = 32830.8
$ws.Range("L63").Value = 32830.8
$ws.Range("N63").Value = -34078.8

# WVR row 65
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 7958.1665
$ws.Range("I65").Value = 2449.5
$ws.Range("J65").Value = 10712.5
$ws.Range("K65").Value = 12247.5
$ws.Range("L65").Value = 53562.5
$ws.Range("M65").Value = -9127.5
$ws.Range("N65").Value = -59802.5

# WVR row 66
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H66").Value = 32830.8
$ws.Range("J66").Value = 32830.8
$ws.Range("L66").Value = 98492.40000000001
$ws.Range("N66").Value = -104732.4

# WVR row 69
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()

# WVR row 72
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()

Write-Output "Edits applied successfully"